# Update gh-pages data: refresh sales/vote counts (column F) across sheets
# 展览 (Exhibition), 演出 (Performance), 本地生活 (Local Life), 全部类型 (All Types)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 777
$ws.Range("F3").Value = 14388
$ws.Range("F4").Value = 14545
$ws.Range("F5").Value = 1369
$ws.Range("F7").Value = 5939
$ws.Range("F8").Value = 993
$ws.Range("F12").Value = 202
$ws.Range("F14").Value = 455
$ws.Range("F16").Value = 1228
$ws.Range("F18").Value = 920
$ws.Range("F20").Value = 2302
$ws.Range("F21").Value = 574
$ws.Range("F22").Value = 831
$ws.Range("F23").Value = 3395
$ws.Range("F25").Value = 320
$ws.Range("F26").Value = 2456
$ws.Range("F27").Value = 612
$ws.Range("F28").Value = 119
$ws.Range("F29").Value = 1340
$ws.Range("F30").Value = 1828
$ws.Range("F31").Value = 1086
$ws.Range("F32").Value = 1447
$ws.Range("F35").Value = 4983
$ws.Range("F36").Value = 4941
$ws.Range("F37").Value = 311
$ws.Range("F39").Value = 687
$ws.Range("F40").Value = 696
$ws.Range("F41").Value = 3319
$ws.Range("F45").Value = 119
$ws.Range("F48").Value = 627
$ws.Range("F49").Value = 306
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 128
$ws.Range("F15").Value = 21
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 7699
$ws.Range("F3").Value = 259
$ws.Range("F4").Value = 894
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 7699
$ws.Range("F3").Value = 777
$ws.Range("F4").Value = 259
$ws.Range("F5").Value = 894
$ws.Range("F7").Value = 14389
$ws.Range("F8").Value = 14546
$ws.Range("F9").Value = 1369
$ws.Range("F11").Value = 5939
$ws.Range("F12").Value = 993
$ws.Range("F13").Value = 128
$ws.Range("F17").Value = 455
$ws.Range("F19").Value = 831
$ws.Range("F20").Value = 3395
$ws.Range("F21").Value = 320
$ws.Range("F22").Value = 2456
$ws.Range("F23").Value = 612
$ws.Range("F24").Value = 119
$ws.Range("F25").Value = 1828
$ws.Range("F28").Value = 21
$ws.Range("F31").Value = 1086
$ws.Range("F32").Value = 1447
$ws.Range("F35").Value = 4983
$ws.Range("F36").Value = 4941
$ws.Range("F37").Value = 311
$ws.Range("F38").Value = 687
$ws.Range("F39").Value = 3319
$ws.Range("F42").Value = 119
$ws.Range("F45").Value = 627
$ws.Range("F46").Value = 306
